# Update dados pesquisadores fioce - atualização 10.2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# dias_defasagem (column E) updates for rows 2-30 (mostly +5 days since last extraction)
$ws.Range("E2").Value = 512
$ws.Range("E3").Value = 116
$ws.Range("E4").Value = 589
$ws.Range("E5").Value = 54
$ws.Range("E6").Value = 81
$ws.Range("E7").Value = 414
$ws.Range("E8").Value = 57
$ws.Range("E9").Value = 118
$ws.Range("E10").Value = 211
$ws.Range("E11").Value = 175
$ws.Range("E12").Value = 153
$ws.Range("E13").Value = 89
$ws.Range("E14").Value = 168
$ws.Range("E15").Value = 95
$ws.Range("E16").Value = 6
$ws.Range("E17").Value = 62
$ws.Range("E18").Value = 21
$ws.Range("E19").Value = 19
$ws.Range("E20").Value = 111
$ws.Range("E21").Value = 238
$ws.Range("E22").Value = 40
$ws.Range("E23").Value = 24
$ws.Range("E24").Value = 53
$ws.Range("E25").Value = 153
$ws.Range("E26").Value = 11
$ws.Range("E27").Value = 48
$ws.Range("E28").Value = 103
$ws.Range("E29").Value = 15
$ws.Range("E30").Value = 278

# Row 31 (Marlos de Medeiros Chaves): ultima_atualizacao, dias_defasagem, qte_artigos_periodicos
$ws.Range("D31").Value = "25/10/2024"
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 8

# dias_defasagem continued for rows 32-34
$ws.Range("E32").Value = 238
$ws.Range("E33").Value = 645
$ws.Range("E34").Value = 183

# Row 35 (Roberto Nicolete): ultima_atualizacao, dias_defasagem, qte_artigos_periodicos
$ws.Range("D35").Value = "25/10/2024"
$ws.Range("E35").Value = 4
$ws.Range("F35").Value = 67

# dias_defasagem continued for rows 36-41
$ws.Range("E36").Value = 42
$ws.Range("E37").Value = 44
$ws.Range("E38").Value = 25
$ws.Range("E39").Value = 62
$ws.Range("E40").Value = 174
$ws.Range("E41").Value = 102
